$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Fill in the "Beat Vegas?" (column G) values for the existing
# games on 1/24/2021 (rows 118-124) that were left blank before.
# ---------------------------------------------------------------
$ws.Range("G118").Value = "Yes"
$ws.Range("G119").Value = "No"
$ws.Range("G120").Value = "No"
$ws.Range("G121").Value = "Yes"
$ws.Range("G122").Value = "Yes"
$ws.Range("G123").Value = "No"
$ws.Range("G124").Value = "No"

# ---------------------------------------------------------------
# Append the new games for 1/25/2021 (rows 125-133).
# ---------------------------------------------------------------
$newGames = @(
    @(44221, "IND", "TOR", -2,    0,     -2),
    @(44221, "ORL", "CHO", -1.5,  5.4,   -6.9),
    @(44221, "DET", "PHI", 5,     12.8,  -7.8),
    @(44221, "BRK", "MIA", -7.5,  -7.2,  -0.3),
    @(44221, "CLE", "LAL", 12,    12.4,  -0.4),
    @(44221, "DAL", "DEN", 2,     3.9,   -1.9),
    @(44221, "CHI", "BOS", 3.5,   -8.5,  12),
    @(44221, "GSW", "MIN", -8.5,  -6.3,  -2.2),
    @(44221, "POR", "OKC", -5,    -1.8,  -3.2)
)

$row = 125
foreach ($game in $newGames) {
    $ws.Cells.Item($row, 1).Value = $game[0]
    $ws.Cells.Item($row, 1).NumberFormat = "yyyy\-mm\-dd"
    $ws.Cells.Item($row, 2).Value = $game[1]
    $ws.Cells.Item($row, 3).Value = $game[2]
    $ws.Cells.Item($row, 4).Value = $game[3]
    $ws.Cells.Item($row, 5).Value = $game[4]
    $ws.Cells.Item($row, 6).Value = $game[5]
    $row++
}

# ---------------------------------------------------------------
# Update the view to show the newly added rows, matching the
# selection/scroll state saved by the original author.
# ---------------------------------------------------------------
[void]$ws.Range("A121").Select()
$excel.ActiveWindow.ScrollRow = 121
[void]$ws.Range("E135").Select()
